# 1013-MS-EI-DB-SAR-REC-NON-RNI-CTPD-SAR-MD-TR-1-LateRepayment-Loanproduct.xlsx
# The product name string had a stray space ("...-Late Repayment") which is
# corrected to "...-LateRepayment" on both the input and output sheets, and
# the ProductLoanInput sheet is made the active sheet/selection instead of
# ProductLoanOutput.

$wb = $excel.ActiveWorkbook

$wsInput  = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

$newProductName = "1013-MS-EI-DB-SAR-REC-NON-RNI-CTPD-SAR-MD-TR-1-LateRepayment"

# Fix the product name text on both sheets (remove stray space before "Repayment")
$wsInput.Range("B1").Value  = $newProductName
$wsOutput.Range("B1").Value = $newProductName

# Update the selection on the output sheet to B1 (no longer the active sheet)
$wsOutput.Range("B1").Select()

# Make ProductLoanInput the active sheet and select B1 there (scrolled to top)
$wsInput.Activate()
$wsInput.Range("B1").Select()
